$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date (45233 -> 45243) for every
# data row (rows 2 through 97). Update the whole block in one go.
$ws.Range("C2:C97").Value = 45243
